$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '29.382.63'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  +0.21%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.842.75'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  -0.01%  '

$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  +0.08%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '239.39'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  -0.27%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.6276'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  -0.05%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.001'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  +0.05%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.07402'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  -0.57%  '

$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  -0.16%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '24.89'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  +2.12%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07716'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  -0.21%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.838.11'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  -0.27%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '4.970'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  -0.15%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.6741'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  -0.70%  '

$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  -1.83%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '81.84'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  -0.06%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '6.265'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  +1.51%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '29.387.88'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  +0.07%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '234.34'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  +2.89%  '

$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  +0.30%  '

$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  +0.07%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '7.307'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  -2.49%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '1.003'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  +0.23%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '157.64'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  -0.93%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '8.504'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  +0.27%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.1346'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  -1.61%  '

$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  -0.91%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.07261'
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  +11.64%  '

$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  +4.54%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.480'
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  -0.01%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '4.037'
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  -1.18%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.039'
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  -1.02%  '

$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  -0.64%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.146'
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  +0.56%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.6961'
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  +0.35%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.574'

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.01831'
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  +0.03%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.804'
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  -1.03%  '

$ws.Range("B39").Value = 'Maker'
$ws.Range("C39").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.233.70'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  -1.81%  '

$ws.Range("B40").Value = 'FraxShare'
$ws.Range("C40").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '6.800'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  +1.08%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.9489'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  +2.46%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.001'

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.994.61'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  -0.65%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '101.05'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  -0.12%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '65.28'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  -0.84%  '

$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  +0.43%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.703'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  -1.31%  '

$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  -1.17%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '8.896'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  -0.99%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.3899'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  -0.64%  '

$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  -2.16%  '
